$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.513731718063354
$ws.Range("B1").Value = 1.613121390342712
$ws.Range("C1").Value = 1.809494256973267
$ws.Range("D1").Value = 2.816995859146118
$ws.Range("E1").Value = 4.123106002807617
